$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing annotation text in G12 ("Finished upto lesson 27" -> "Finished upto lesson# 27")
$ws.Range("G12").Value = "Finished upto lesson# 27"

# Fill in the new timesheet row 13 (date, start time, end time)
$ws.Range("B13").Value = 44825
$ws.Range("C13").Value = 0.90625
$ws.Range("D13").Value = 0.95833333333333337

# Compute the elapsed time for the new row, same as the other rows in the table
$ws.Range("E13").Formula = "=D13-C13"

# Add the new progress note for row 13, matching the formatting used for the
# note in G12 (copy format only so the text we just set is preserved)
$ws.Range("G13").Value = "Finished upto lesson# 29"
$ws.Range("G12").Copy()
$ws.Range("G13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Move the active selection to reflect where the user ended up editing
$ws.Range("G14").Select()
